$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Control 1
$ws.Range("D2").Value = 0.9999986599940173
$ws.Range("E2").Value = 0.9999986599940173

# Row 3 - Control 31
$ws.Range("D3").Value = 0.9999999983753163
$ws.Range("E3").Value = 0.9999999983753163

# Row 4 - Control 14
$ws.Range("D4").Value = 0.9316902151288704
$ws.Range("E4").Value = 0.9316902151288704

# Row 5 - Control 19
$ws.Range("D5").Value = 0.000000004333885456476208
$ws.Range("E5").Value = 0.000000004333885456476208

# Row 6 - MDD 27
$ws.Range("D6").Value = 0.007776851039568014
$ws.Range("E6").Value = 0.007776851039568014

# Row 7 - MDD 47
$ws.Range("D7").Value = 0.9999991136584325
$ws.Range("E7").Value = 0.000000886341567474247

# Row 8 - MDD 13
$ws.Range("D8").Value = 0.0000002541333987880357
$ws.Range("E8").Value = 0.9999997458666012

# Row 9 - MDD 25
$ws.Range("D9").Value = 0.9567689275690073
$ws.Range("E9").Value = 0.04323107243099267

# Row 11 - MDD 5
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.03082207722633906
$ws.Range("E11").Value = 0.969177922773661
$ws.Range("F11").Value = 5.516141891479492
$ws.Range("G11").Value = 0.5
